$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated column widths (from diff). The engine quantizes ColumnWidth to a
# display-pixel grid (MDW=7 => 1/6-character steps) before it is stored back
# as the OOXML `width` (in characters), so we pick the ColumnWidth input
# whose quantized, stored result is closest to the target stored width
# (target 15.7109375 -> stored 15.666666..., target 16.42578125 -> stored 16.5).
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666

$colA = @(
-0.31112539381262394,-0.20818866504045452,-0.14429667797497103,-0.13594953726081016,-0.1322842600741554,-0.03075611900198716,-0.020465560411039974,-0.010406153339351576,-0.008328440261988312,-0.00627360289652934,-0.003273322363856579,0.00022960687714235206,0.003748781132436818,0.011753247918931287,0.012757392222026098,-0.006033846931785103,-0.004003096305788212,-0.044643076472219434,-0.012090857353391993,-0.008015931530204412,-0.004005563956327229,-0.04570619882239946,-0.04049503324510173,-0.020098067566756583,-0.09723559439944296,-0.09461287832594145,-0.09195425279488978,-0.08900702143887251,-0.08135874876256732,-0.021171358891272085,-0.014022014230643975,-0.004000733624154762
)
$colB = @(
0.31077622830142104,0.207467366934031,0.14394953714330327,0.13528425999906624,0.12999171330923165,0.030465560234491207,0.020406153159997942,0.010328440163239083,0.008273602795931367,0.006273322252690505,0.0032703930063808784,-0.00024878124892202536,-0.0037532480817867864,-0.01175739231126105,-0.012765263253087689,0.006003096206049996,0.003999999879072291,0.04459424158615022,0.012015931480176878,0.008005563905811641,0.003999999949055422,0.04549503316474102,0.040098067323480535,0.019999999753369302,0.09711287822617365,0.09445425269009178,0.0910070213146037,0.08835874857008363,0.08117135814085863,0.02102201402438908,0.014000733386069442,0.003999999823859568
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
}
